# Insert two new data rows right before the current row 466.
# This shifts the existing rows 466-567 down to 468-569 (matching the
# target dimension A1:R569), and leaves room for the two new rows of
# data (new row 466 and new row 467).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("466:467").Insert()

# --- New row 466 ---
$ws.Cells.Item(466, 1).Value  = 3
$ws.Cells.Item(466, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(466, 3).Value  = "Coquimbo"
$ws.Cells.Item(466, 4).Value  = 44889
$ws.Cells.Item(466, 5).Value  = 5
$ws.Cells.Item(466, 6).Value  = 100112032
$ws.Cells.Item(466, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(466, 8).Value  = "Sin especificar"
$ws.Cells.Item(466, 9).Value  = "Primera"
$ws.Cells.Item(466, 10).Value = 90
$ws.Cells.Item(466, 11).Value = 4000
$ws.Cells.Item(466, 12).Value = 4000
$ws.Cells.Item(466, 13).Value = 4000
$ws.Cells.Item(466, 14).Value = "$/caja 36 unidades"
$ws.Cells.Item(466, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(466, 16).Value = 111
$ws.Cells.Item(466, 17).Value = 36
$ws.Cells.Item(466, 18).Value = "Hortaliza"

# --- New row 467 ---
$ws.Cells.Item(467, 1).Value  = 3
$ws.Cells.Item(467, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(467, 3).Value  = "Coquimbo"
$ws.Cells.Item(467, 4).Value  = 44889
$ws.Cells.Item(467, 5).Value  = 5
$ws.Cells.Item(467, 6).Value  = 100112032
$ws.Cells.Item(467, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(467, 8).Value  = "Sin especificar"
$ws.Cells.Item(467, 9).Value  = "Primera"
$ws.Cells.Item(467, 10).Value = 190
$ws.Cells.Item(467, 11).Value = 7500
$ws.Cells.Item(467, 12).Value = 8000
$ws.Cells.Item(467, 13).Value = 7747
$ws.Cells.Item(467, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(467, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(467, 16).Value = 111
$ws.Cells.Item(467, 17).Value = 70
$ws.Cells.Item(467, 18).Value = "Hortaliza"
